# Apply cryptos list update (Wed Jul 19 15:24:19 UTC 2023, GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Text)
    # Force literal text into a cell even when it looks numeric (e.g. '0.9995',
    # '1.000'), while preserving the cell's original (unset) style.
    $origStyle = $Cell.Style
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = $origStyle
}

$ws.Range('D2').Value = '29.871.48'
$ws.Range('E2').Value = '  -0.11%  '
$ws.Range('D3').Value = '1.900.57'
$ws.Range('E3').Value = '  +0.18%  '
Set-TextValue $ws.Range('D4') '0.9995'
$ws.Range('E4').Value = '  -0.10%  '
Set-TextValue $ws.Range('D5') '0.8033'
$ws.Range('E5').Value = '  +6.27%  '
Set-TextValue $ws.Range('D6') '241.15'
$ws.Range('E6').Value = '  +0.44%  '
Set-TextValue $ws.Range('D7') '1.000'
Set-TextValue $ws.Range('D8') '0.3125'
$ws.Range('E8').Value = '  +2.69%  '
Set-TextValue $ws.Range('D9') '26.13'
$ws.Range('E9').Value = '  +2.91%  '
Set-TextValue $ws.Range('D10') '0.06878'
$ws.Range('E10').Value = '  +0.60%  '
Set-TextValue $ws.Range('D11') '0.07974'
$ws.Range('E11').Value = '  -0.07%  '
$ws.Range('D12').Value = '1.903.00'
$ws.Range('E12').Value = '  +0.11%  '
Set-TextValue $ws.Range('D13') '0.7362'
$ws.Range('E13').Value = '  -1.61%  '
Set-TextValue $ws.Range('D14') '5.162'
$ws.Range('E14').Value = '  -0.88%  '
Set-TextValue $ws.Range('D15') '92.33'
$ws.Range('E15').Value = '  +1.20%  '
$ws.Range('D16').Value = '29.867.40'
$ws.Range('E16').Value = '  -0.17%  '
Set-TextValue $ws.Range('D17') '13.88'
$ws.Range('E17').Value = '  -0.46%  '
Set-TextValue $ws.Range('D18') '5.835'
Set-TextValue $ws.Range('D19') '244.51'
$ws.Range('E19').Value = '  +0.63%  '
Set-TextValue $ws.Range('D20') '0.000007698'
$ws.Range('E20').Value = '  -0.34%  '
$ws.Range('E21').Value = '  +0.01%  '
$ws.Range('D22').Value = '2.151.30'
$ws.Range('E22').Value = '  -0.60%  '
Set-TextValue $ws.Range('D23') '0.9994'
$ws.Range('E23').Value = '  -0.14%  '
Set-TextValue $ws.Range('D24') '6.844'
$ws.Range('E24').Value = '  -1.56%  '
Set-TextValue $ws.Range('D25') '166.89'
$ws.Range('E25').Value = '  +0.61%  '
Set-TextValue $ws.Range('D26') '9.171'
$ws.Range('E26').Value = '  -0.69%  '
Set-TextValue $ws.Range('D27') '0.1423'
$ws.Range('E27').Value = '  +9.06%  '
Set-TextValue $ws.Range('D28') '18.81'
$ws.Range('E28').Value = '  +0.45%  '
Set-TextValue $ws.Range('D29') '2.026'
$ws.Range('E29').Value = '  +0.54%  '
Set-TextValue $ws.Range('D30') '1.355'
$ws.Range('E30').Value = '  -4.06%  '
Set-TextValue $ws.Range('D31') '1.513'
$ws.Range('E31').Value = '  -0.39%  '
Set-TextValue $ws.Range('D32') '4.281'
$ws.Range('E32').Value = '  +0.07%  '
Set-TextValue $ws.Range('D33') '0.05559'
$ws.Range('E33').Value = '  +3.45%  '
Set-TextValue $ws.Range('D34') '4.060'
$ws.Range('E34').Value = '  +0.88%  '
Set-TextValue $ws.Range('D35') '1.254'
$ws.Range('E35').Value = '  +0.07%  '
Set-TextValue $ws.Range('D36') '0.7246'
$ws.Range('E36').Value = '  -0.14%  '
$ws.Range('E37').Value = '  +0.08%  '
$ws.Range('E38').Value = '  +0.14%  '
Set-TextValue $ws.Range('D39') '2.780'
$ws.Range('E39').Value = '  -0.39%  '
Set-TextValue $ws.Range('D40') '0.4387'
$ws.Range('E40').Value = '  -0.37%  '
Set-TextValue $ws.Range('D41') '5.995'
$ws.Range('E41').Value = '  -2.79%  '
Set-TextValue $ws.Range('D42') '71.89'
$ws.Range('E42').Value = '  -0.49%  '
$ws.Range('E43').Value = '  -0.03%  '
Set-TextValue $ws.Range('D44') '0.8334'
$ws.Range('E44').Value = '  +0.84%  '
Set-TextValue $ws.Range('D45') '1.851'
$ws.Range('E45').Value = '  -3.01%  '
Set-TextValue $ws.Range('D46') '100.45'
$ws.Range('E46').Value = '  -0.64%  '
Set-TextValue $ws.Range('D47') '7.546'
$ws.Range('E47').Value = '  -0.31%  '
Set-TextValue $ws.Range('D48') '9.718'
$ws.Range('E48').Value = '  -0.79%  '
$ws.Range('B49').Value = 'RocketPoolETH'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D49').Value = '2.060.60'
$ws.Range('E49').Value = '  +0.09%  '
$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws.Range('D50') '975.91'
$ws.Range('E50').Value = '  +7.36%  '
Set-TextValue $ws.Range('D51') '36.09'
$ws.Range('E51').Value = '  -0.51%  '
